$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Is Active") currently stores text strings "True"/"False".
# Push up set value: change these to native boolean values (TRUE/FALSE)
# so the cells are stored as boolean-typed cells, not shared strings.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
